$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1999.4445
$ws.Range("I98").Value = 2070.7144
$ws.Range("K98").Value = 2070.7144
$ws.Range("M98").Value = -572.7143999999998

$ws.Range("H122").Value = 1999.4445
$ws.Range("I122").Value = 2070.7144
$ws.Range("K122").Value = 6212.1432
$ws.Range("M122").Value = -3762.1432

$ws.Range("H135").Value = 13890541
$ws.Range("I135").Value = 1254.8572
$ws.Range("J135").Value = 62503040
$ws.Range("K135").Value = 11293.7148
$ws.Range("L135").Value = 562527360
$ws.Range("M135").Value = -8758.7148
$ws.Range("N135").Value = -562532430

$ws.Range("H137").Value = 6250961.5
$ws.Range("I137").Value = 928.2593000000001
$ws.Range("J137").Value = 40001140
$ws.Range("K137").Value = 2784.7779
$ws.Range("L137").Value = 120003420
$ws.Range("M137").Value = -234.7779
$ws.Range("N137").Value = -120008520

$ws.Range("H139").Value = 44288.125
$ws.Range("I139").Value = 20709
$ws.Range("J139").Value = 58435.6
$ws.Range("K139").Value = 20709
$ws.Range("L139").Value = 58435.6
$ws.Range("M139").Value = -15569
$ws.Range("N139").Value = -68715.60000000001

$ws.Range("H140").Value = 78860
$ws.Range("J140").Value = 78860
$ws.Range("L140").Value = 78860
$ws.Range("N140").Value = -89220

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8108.536
$ws.Range("I32").Value = 8121.9204
$ws.Range("J32").Value = 7968
$ws.Range("K32").Value = 8121.9204
$ws.Range("L32").Value = 7968
$ws.Range("M32").Value = -7834.9204
$ws.Range("N32").Value = -8542

$ws.Range("H61").Value = 13890370
$ws.Range("I61").Value = 17242886
$ws.Range("J61").Value = 1371.4286
$ws.Range("K61").Value = 17242886
$ws.Range("L61").Value = 1371.4286
$ws.Range("M61").Value = -17242674
$ws.Range("N61").Value = -1795.4286

$ws.Range("H74").Value = 18522678
$ws.Range("I74").Value = 38464104
$ws.Range("J74").Value = 5639.2856
$ws.Range("K74").Value = 38464104
$ws.Range("L74").Value = 5639.2856
$ws.Range("M74").Value = -38463230
$ws.Range("N74").Value = -7387.2856

$ws.Range("H77").Value = 18522678
$ws.Range("I77").Value = 38464104
$ws.Range("J77").Value = 5639.2856
$ws.Range("K77").Value = 192320520
$ws.Range("L77").Value = 28196.428
$ws.Range("M77").Value = -192316152
$ws.Range("N77").Value = -36932.428

$ws.Range("H132").Value = 5816238.5
$ws.Range("I132").Value = 8066345.5
$ws.Range("J132").Value = 3461.9167
$ws.Range("K132").Value = 24199036.5
$ws.Range("L132").Value = 10385.7501
$ws.Range("M132").Value = -24196506.5
$ws.Range("N132").Value = -15445.7501

$ws.Range("H136").Value = 13890370
$ws.Range("I136").Value = 17242886
$ws.Range("J136").Value = 1371.4286
$ws.Range("K136").Value = 51728658
$ws.Range("L136").Value = 4114.2858
$ws.Range("M136").Value = -51726108
$ws.Range("N136").Value = -9214.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3775.1562
$ws.Range("I134").Value = 2656.5715
$ws.Range("J134").Value = 5910.636
$ws.Range("K134").Value = 7969.7145
$ws.Range("L134").Value = 17731.908
$ws.Range("M134").Value = -5434.7145
$ws.Range("N134").Value = -22801.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9010777
$ws.Range("I31").Value = 1840.5938
$ws.Range("J31").Value = 66667972
$ws.Range("K31").Value = 1840.5938
$ws.Range("L31").Value = 66667972
$ws.Range("M31").Value = -1545.5938
$ws.Range("N31").Value = -66668562

$ws.Range("H34").Value = 9010777
$ws.Range("I34").Value = 1840.5938
$ws.Range("J34").Value = 66667972
$ws.Range("K34").Value = 1840.5938
$ws.Range("L34").Value = 66667972
$ws.Range("M34").Value = -1638.5938
$ws.Range("N34").Value = -66668376

$ws.Range("H58").Value = 1828.5294
$ws.Range("I58").Value = 805.2941
$ws.Range("J58").Value = 2851.7646
$ws.Range("K58").Value = 805.2941
$ws.Range("L58").Value = 2851.7646
$ws.Range("M58").Value = -602.2941
$ws.Range("N58").Value = -3257.7646

$ws.Range("H132").Value = 17859800
$ws.Range("I132").Value = 22729198
$ws.Range("J132").Value = 5337.6665
$ws.Range("K132").Value = 68187594
$ws.Range("L132").Value = 16012.9995
$ws.Range("M132").Value = -68185064
$ws.Range("N132").Value = -21072.9995

$ws.Range("H134").Value = 1833605.2
$ws.Range("I134").Value = 2395
$ws.Range("J134").Value = 11905262
$ws.Range("K134").Value = 7185
$ws.Range("L134").Value = 35715786
$ws.Range("M134").Value = -4650
$ws.Range("N134").Value = -35720856

$ws.Range("H136").Value = 1828.5294
$ws.Range("I136").Value = 805.2941
$ws.Range("J136").Value = 2851.7646
$ws.Range("K136").Value = 2415.8823
$ws.Range("L136").Value = 8555.293799999999
$ws.Range("M136").Value = 134.1177000000002
$ws.Range("N136").Value = -13655.2938

$ws.Range("H140").Value = 33493.7
$ws.Range("J140").Value = 33493.7
$ws.Range("L140").Value = 33493.7
$ws.Range("N140").Value = -43853.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1000
$ws.Range("J117").Value = 1000
$ws.Range("L117").Value = 3000
$ws.Range("N117").Value = -9884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3447.7844
$ws.Range("I132").Value = 2396.5588
$ws.Range("J132").Value = 5550.2354
$ws.Range("K132").Value = 7189.676399999999
$ws.Range("L132").Value = 16650.7062
$ws.Range("M132").Value = -4659.676399999999
$ws.Range("N132").Value = -21710.7062

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5419.7393
$ws.Range("I40").Value = 4603.375
$ws.Range("J40").Value = 7285.7144
$ws.Range("K40").Value = 4603.375
$ws.Range("L40").Value = 7285.7144
$ws.Range("M40").Value = -4467.375
$ws.Range("N40").Value = -7557.7144

$ws.Range("H132").Value = 6678.0894
$ws.Range("I132").Value = 4472.5864
$ws.Range("J132").Value = 9046.963
$ws.Range("K132").Value = 13417.7592
$ws.Range("L132").Value = 27140.889
$ws.Range("M132").Value = -10887.7592
$ws.Range("N132").Value = -32200.889

$ws.Range("H136").Value = 14711509
$ws.Range("I136").Value = 21741282
$ws.Range("J136").Value = 12891.637
$ws.Range("K136").Value = 65223846
$ws.Range("L136").Value = 38674.911
$ws.Range("M136").Value = -65221296
$ws.Range("N136").Value = -43774.911

$ws.Range("H139").Value = 53389.3
$ws.Range("J139").Value = 54804.777
$ws.Range("L139").Value = 54804.777
$ws.Range("N139").Value = -65084.777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1895.9474
$ws.Range("I132").Value = 1258.12
$ws.Range("J132").Value = 3122.5386
$ws.Range("K132").Value = 3774.36
$ws.Range("L132").Value = 9367.6158
$ws.Range("M132").Value = -1244.36
$ws.Range("N132").Value = -14427.6158

$ws.Range("H136").Value = 1135.0416
$ws.Range("I136").Value = 1182.9048
$ws.Range("J136").Value = 800
$ws.Range("K136").Value = 3548.7144
$ws.Range("L136").Value = 2400
$ws.Range("M136").Value = -998.7143999999998
$ws.Range("N136").Value = -7500
